# "Generate Report for Handoff"
# Adds a new tracked file (b081659a-a9f2-42a7-ad0f-7e0c917c1f7e...) that is
# "Ready for handoff" to the Overview sheet and to the per-locale
# (zh-cn / de-de) detail sheets, each backed by an Excel Table.

$wb = $excel.ActiveWorkbook

$newGuidMd        = "b081659a-a9f2-42a7-ad0f-7e0c917c1f7eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newGuidMdWithDir = "e2e\b081659a-a9f2-42a7-ad0f-7e0c917c1f7eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$readyStatus      = "Ready for handoff"
$handoffDate      = "2016-09-01 10:32:14"
$zhXlf            = "b081659a-a9f2-42a7-ad0f-7e0c917c1f7eoooooooooooooooooooooooooooooooooooooooo.49b8a2520e9e2df664159de762e2eea0e9536db7.zh-cn.xlf"
$zhXlfDate        = "2016-09-01 10:32:02"
$deXlf            = "b081659a-a9f2-42a7-ad0f-7e0c917c1f7eoooooooooooooooooooooooooooooooooooooooo.49b8a2520e9e2df664159de762e2eea0e9536db7.de-de.xlf"

$newMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c74cb38ba9f6a7473ac8ced7be8c4fb87144fa6/e2e/b081659a-a9f2-42a7-ad0f-7e0c917c1f7eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newGuidMd
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $handoffDate
$wsOverview.Range("G3").NumberFormatLocal = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdUrl, "", "", $newGuidMdWithDir) | Out-Null
$wsOverview.Range("B3").Font.Underline = 2
$wsOverview.Range("B3").Font.Color = 15570276

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyStatus
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhXlfDate
$wsZh.Range("H3").NumberFormatLocal = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormatLocal = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newMdUrl, "", "", $newGuidMd) | Out-Null
$wsZh.Range("A3").Font.Underline = 2
$wsZh.Range("A3").Font.Color = 15570276

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyStatus
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $handoffDate
$wsDe.Range("H3").NumberFormatLocal = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormatLocal = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newMdUrl, "", "", $newGuidMd) | Out-Null
$wsDe.Range("A3").Font.Underline = 2
$wsDe.Range("A3").Font.Color = 15570276

Write-Host "Report generated for handoff."
